$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("newpage1")
$ws.Activate()

# Set C3 to the same email address as C2 ("Shilpi Mow"'s row), with a mailto hyperlink,
# mirroring the existing hyperlinked email already present in C2.
$ws.Range("C3").Value = "shilpimou@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:shilpimou@gmail.com")
$ws.Range("C3").Style = "Hyperlink"

# Update the selection to match the new authored state
$ws.Range("C2:C3").Select()
